# Generate Report for Handback
# Update the handoff/handback timestamp strings on the zh-cn and de-de
# report sheets to reflect the newly generated report run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 12:54:12"
$wsZhCn.Range("H2").Value = "2016-03-20 12:54:31"

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 12:54:15"
$wsDeDe.Range("H2").Value = "2016-03-20 12:54:37"
